# Bugfixed evaluation and simulated rt_data for components.
#
# The YoY forecast-vector sheet gets one more observation inserted at the
# front (a new oldest forecast, 2007->2008, at new row 2) and one more
# appended at the end (a new newest forecast, 2025->2026, at new row 53).
# Every existing observation shifts down one row, and because the
# evaluation bug fix also changed how y_0_forecast / y_1_forecast are
# simulated, the C and E columns are recomputed for every single row (A/B/D
# - the actual date and the y_0/y_1 calendar years - stay the same values,
# just shifted down one row). Simplest + most robust way to land all of
# that is to just (re)write the full A2:E53 block with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 1.75539628881467
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 1.327368416067398

$ws.Range("A3").Value = 39583
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = 1.392321641630434
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = 1.710071460977503

$ws.Range("A4").Value = 39765
$ws.Range("B4").Value = 2008
$ws.Range("C4").Value = 2.213911448916162
$ws.Range("D4").Value = 2009
$ws.Range("E4").Value = 2.649257112350067

$ws.Range("A5").Value = 39948
$ws.Range("B5").Value = 2009
$ws.Range("C5").Value = 1.004409005705997
$ws.Range("D5").Value = 2010
$ws.Range("E5").Value = 1.642433761320072

$ws.Range("A6").Value = 40130
$ws.Range("B6").Value = 2009
$ws.Range("C6").Value = 2.533533936850563
$ws.Range("D6").Value = 2010
$ws.Range("E6").Value = 1.815660192323709

$ws.Range("A7").Value = 40310
$ws.Range("B7").Value = 2010
$ws.Range("C7").Value = 2.418114148635109
$ws.Range("D7").Value = 2011
$ws.Range("E7").Value = 2.828066716168021

$ws.Range("A8").Value = 40494
$ws.Range("B8").Value = 2010
$ws.Range("C8").Value = 2.088987486264915
$ws.Range("D8").Value = 2011
$ws.Range("E8").Value = 2.332261646026201

$ws.Range("A9").Value = 40676
$ws.Range("B9").Value = 2011
$ws.Range("C9").Value = 1.89159218653383
$ws.Range("D9").Value = 2012
$ws.Range("E9").Value = 2.544631191216329

$ws.Range("A10").Value = 40862
$ws.Range("B10").Value = 2011
$ws.Range("C10").Value = 1.212544822741002
$ws.Range("D10").Value = 2012
$ws.Range("E10").Value = 1.839804681163293

$ws.Range("A11").Value = 41044
$ws.Range("B11").Value = 2012
$ws.Range("C11").Value = 1.554977796875501
$ws.Range("D11").Value = 2013
$ws.Range("E11").Value = 1.312870290004287

$ws.Range("A12").Value = 41228
$ws.Range("B12").Value = 2012
$ws.Range("C12").Value = 1.196776590518644
$ws.Range("D12").Value = 2013
$ws.Range("E12").Value = 0.670590452940556

$ws.Range("A13").Value = 41409
$ws.Range("B13").Value = 2013
$ws.Range("C13").Value = 0.6180254938795482
$ws.Range("D13").Value = 2014
$ws.Range("E13").Value = 0.7749619016293785

$ws.Range("A14").Value = 41592
$ws.Range("B14").Value = 2013
$ws.Range("C14").Value = 0.4712609263772594
$ws.Range("D14").Value = 2014
$ws.Range("E14").Value = 0.8520644823059476

$ws.Range("A15").Value = 41774
$ws.Range("B15").Value = 2014
$ws.Range("C15").Value = 0.481899667566732
$ws.Range("D15").Value = 2015
$ws.Range("E15").Value = 0.7487574275252262

$ws.Range("A16").Value = 41957
$ws.Range("B16").Value = 2014
$ws.Range("C16").Value = 0.8783377572271434
$ws.Range("D16").Value = 2015
$ws.Range("E16").Value = 1.474590898715178

$ws.Range("A17").Value = 42137
$ws.Range("B17").Value = 2015
$ws.Range("C17").Value = 1.905862317202089
$ws.Range("D17").Value = 2016
$ws.Range("E17").Value = 1.389591155234515

$ws.Range("A18").Value = 42321
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = 2.29066283401107
$ws.Range("D18").Value = 2016
$ws.Range("E18").Value = 2.597902967862775

$ws.Range("A19").Value = 42503
$ws.Range("B19").Value = 2016
$ws.Range("C19").Value = 2.671046044496239
$ws.Range("D19").Value = 2017
$ws.Range("E19").Value = 2.125743999456575

$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = 2016
$ws.Range("C20").Value = 4.109890522944348
$ws.Range("D20").Value = 2017
$ws.Range("E20").Value = 3.628019428949036

$ws.Range("A21").Value = 42867
$ws.Range("B21").Value = 2017
$ws.Range("C21").Value = 1.917627847674064
$ws.Range("D21").Value = 2018
$ws.Range("E21").Value = 2.694711744616662

$ws.Range("A22").Value = 43053
$ws.Range("B22").Value = 2017
$ws.Range("C22").Value = 1.336316831462692
$ws.Range("D22").Value = 2018
$ws.Range("E22").Value = 1.626630409005325

$ws.Range("A23").Value = 43145
$ws.Range("B23").Value = 2018
$ws.Range("C23").Value = 2.121911365876805
$ws.Range("D23").Value = 2019
$ws.Range("E23").Value = 2.706722015217466

$ws.Range("A24").Value = 43235
$ws.Range("B24").Value = 2018
$ws.Range("C24").Value = 1.119562422009102
$ws.Range("D24").Value = 2019
$ws.Range("E24").Value = 1.831617848540201

$ws.Range("A25").Value = 43326
$ws.Range("B25").Value = 2018
$ws.Range("C25").Value = 1.33730574578026
$ws.Range("D25").Value = 2019
$ws.Range("E25").Value = 1.730186041121162

$ws.Range("A26").Value = 43418
$ws.Range("B26").Value = 2018
$ws.Range("C26").Value = 1.197912858979611
$ws.Range("D26").Value = 2019
$ws.Range("E26").Value = 1.216371234267344

$ws.Range("A27").Value = 43510
$ws.Range("B27").Value = 2019
$ws.Range("C27").Value = 2.252616573494293
$ws.Range("D27").Value = 2020
$ws.Range("E27").Value = 1.465076104875918

$ws.Range("A28").Value = 43600
$ws.Range("B28").Value = 2019
$ws.Range("C28").Value = 1.344920716048192
$ws.Range("D28").Value = 2020
$ws.Range("E28").Value = 1.037735724446631

$ws.Range("A29").Value = 43691
$ws.Range("B29").Value = 2019
$ws.Range("C29").Value = 1.722359355602787
$ws.Range("D29").Value = 2020
$ws.Range("E29").Value = 1.75179450213927

$ws.Range("A30").Value = 43783
$ws.Range("B30").Value = 2019
$ws.Range("C30").Value = 1.727537197898665
$ws.Range("D30").Value = 2020
$ws.Range("E30").Value = 2.164378481800822

$ws.Range("A31").Value = 43875
$ws.Range("B31").Value = 2020
$ws.Range("C31").Value = 2.554068495740247
$ws.Range("D31").Value = 2021
$ws.Range("E31").Value = 2.057677568601401

$ws.Range("A32").Value = 43966
$ws.Range("B32").Value = 2020
$ws.Range("C32").Value = 2.195375580740766
$ws.Range("D32").Value = 2021
$ws.Range("E32").Value = 1.872521508785896

$ws.Range("A33").Value = 44068
$ws.Range("B33").Value = 2020
$ws.Range("C33").Value = 2.980209378995857
$ws.Range("D33").Value = 2021
$ws.Range("E33").Value = 2.822333853751413

$ws.Range("A34").Value = 44159
$ws.Range("B34").Value = 2020
$ws.Range("C34").Value = 3.647228437274408
$ws.Range("D34").Value = 2021
$ws.Range("E34").Value = 3.845906281600109

$ws.Range("A35").Value = 44251
$ws.Range("B35").Value = 2021
$ws.Range("C35").Value = 2.115981176209125
$ws.Range("D35").Value = 2022
$ws.Range("E35").Value = 2.767253381388879

$ws.Range("A36").Value = 44341
$ws.Range("B36").Value = 2021
$ws.Range("C36").Value = 2.542856270410665
$ws.Range("D36").Value = 2022
$ws.Range("E36").Value = 2.961494745505977

$ws.Range("A37").Value = 44432
$ws.Range("B37").Value = 2021
$ws.Range("C37").Value = 1.954146674711188
$ws.Range("D37").Value = 2022
$ws.Range("E37").Value = 1.581547781257497

$ws.Range("A38").Value = 44525
$ws.Range("B38").Value = 2021
$ws.Range("C38").Value = 2.777797690741424
$ws.Range("D38").Value = 2022
$ws.Range("E38").Value = 1.875884305456199

$ws.Range("A39").Value = 44617
$ws.Range("B39").Value = 2022
$ws.Range("C39").Value = 1.450993313666182
$ws.Range("D39").Value = 2023
$ws.Range("E39").Value = 2.17564691785852

$ws.Range("A40").Value = 44706
$ws.Range("B40").Value = 2022
$ws.Range("C40").Value = 0.3979826440748235
$ws.Range("D40").Value = 2023
$ws.Range("E40").Value = 2.008592810942544

$ws.Range("A41").Value = 44798
$ws.Range("B41").Value = 2022
$ws.Range("C41").Value = 2.69102598245059
$ws.Range("D41").Value = 2023
$ws.Range("E41").Value = 4.307673059319161

$ws.Range("A42").Value = 44890
$ws.Range("B42").Value = 2022
$ws.Range("C42").Value = 0.6994919452575576
$ws.Range("D42").Value = 2023
$ws.Range("E42").Value = -0.2388228654152447

$ws.Range("A43").Value = 44981
$ws.Range("B43").Value = 2023
$ws.Range("C43").Value = -0.1329858710789389
$ws.Range("D43").Value = 2024
$ws.Range("E43").Value = 0.4083040303828334

$ws.Range("A44").Value = 45071
$ws.Range("B44").Value = 2023
$ws.Range("C44").Value = -2.604000402888396
$ws.Range("D44").Value = 2024
$ws.Range("E44").Value = -0.08252516517808228

$ws.Range("A45").Value = 45163
$ws.Range("B45").Value = 2023
$ws.Range("C45").Value = -1.669605379075589
$ws.Range("D45").Value = 2024
$ws.Range("E45").Value = -0.3682427893006324

$ws.Range("A46").Value = 45254
$ws.Range("B46").Value = 2023
$ws.Range("C46").Value = -1.432689847121871
$ws.Range("D46").Value = 2024
$ws.Range("E46").Value = -0.7896638887521124

$ws.Range("A47").Value = 45345
$ws.Range("B47").Value = 2024
$ws.Range("C47").Value = 0.6825239311359033
$ws.Range("D47").Value = 2025
$ws.Range("E47").Value = 0.06285237552883238

$ws.Range("A48").Value = 45436
$ws.Range("B48").Value = 2024
$ws.Range("C48").Value = 1.122551915563408
$ws.Range("D48").Value = 2025
$ws.Range("E48").Value = 0.254631175783615

$ws.Range("A49").Value = 45534
$ws.Range("B49").Value = 2024
$ws.Range("C49").Value = 1.780300968358017
$ws.Range("D49").Value = 2025
$ws.Range("E49").Value = 0.8629164812201218

$ws.Range("A50").Value = 45618
$ws.Range("B50").Value = 2024
$ws.Range("C50").Value = 2.033479419175133
$ws.Range("D50").Value = 2025
$ws.Range("E50").Value = 1.424898175306621

$ws.Range("A51").Value = 45713
$ws.Range("B51").Value = 2025
$ws.Range("C51").Value = 2.87910633698536
$ws.Range("D51").Value = 2026
$ws.Range("E51").Value = 1.829048896543739

$ws.Range("A52").Value = 45800
$ws.Range("B52").Value = 2025
$ws.Range("C52").Value = 2.273132718878146
$ws.Range("D52").Value = 2026
$ws.Range("E52").Value = 1.620205313802381

$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 2025
$ws.Range("C53").Value = 2.481068287768839
$ws.Range("D53").Value = 2026
$ws.Range("E53").Value = 2.274625453684709

# Row 53 is brand new (the sheet previously only went to row 52), so it has
# no formatting yet. Column A throughout the table carries the date display
# style (centered, bordered, YYYY-MM-DD HH:MM:SS number format) - copy that
# formatting down from A52 so the new row matches the rest of the column.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "Forecast table rewritten; used range is now $($ws.UsedRange.Address())"
